$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Range("D12").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E12").Value = "['Normal']"

# Row 38
$ws.Range("D38").Value = "[0, 0, 1, 0, 0, 0, 0]"
$ws.Range("E38").Value = "['HardwareFault']"

# Row 39
$ws.Range("D39").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E39").Value = "['Normal', 'HardwareFault']"

# Row 58
$ws.Range("D58").Value = "[0, 0, 0, 1, 0, 0, 0]"
$ws.Range("E58").Value = "['ParamViolation']"

# Row 61
$ws.Range("D61").Value = "[0, 0, 0, 0, 0, 0, 1]"
$ws.Range("E61").Value = "['SoftwareFault']"

# Row 67
$ws.Range("D67").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E67").Value = "['Normal']"

# Row 68
$ws.Range("D68").Value = "[1, 0, 0, 1, 0, 0, 0]"
$ws.Range("E68").Value = "['Normal', 'ParamViolation']"

# Row 71
$ws.Range("D71").Value = "[1, 0, 0, 1, 0, 0, 0]"
$ws.Range("E71").Value = "['Normal', 'ParamViolation']"

# Row 73
$ws.Range("D73").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E73").Value = "['Normal', 'SoftwareFault']"

# Row 75
$ws.Range("D75").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E75").Value = "['Normal', 'SoftwareFault']"

# Row 84
$ws.Range("D84").Value = "[1, 0, 0, 1, 0, 0, 0]"
$ws.Range("E84").Value = "['Normal', 'ParamViolation']"

# Row 88
$ws.Range("D88").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E88").Value = "['Normal', 'HardwareFault']"

# Row 92
$ws.Range("D92").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E92").Value = "['Normal', 'SoftwareFault']"

# Row 107
$ws.Range("D107").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E107").Value = "['Normal', 'SoftwareFault']"
